# Append the new sales row (row 10) to the bottom of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "date" as plain text in this sheet (same as the rest of
# the column), so force Text formatting while entering the value to stop
# Excel from auto-converting the "2025-01-05" string into a date serial
# number, then restore the default "Normal" style so no stray formatting
# is left behind on the new cell.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-01-05"
$ws.Range("A10").Style = "Normal"

$ws.Range("B10").Value = "Phone"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 499.32
$ws.Range("E10").Value = 998.64
